# gridsearch.xlsx update: add NU_F / NU_V columns, add 3 new result rows,
# and resize/relabel the Score column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Move the existing "SCORE" column (D) and "Note" column (E) two
#    columns to the right, to F and G, to make room for the new
#    NU_F / NU_V columns at D and E.  Only copy the cells that actually
#    contain data so we do not materialize stray blank cells.
# ---------------------------------------------------------------------
$ws.Range("E1").Copy($ws.Range("G1"))
$ws.Range("E12:E23").Copy($ws.Range("G12:G23"))
$ws.Range("D1:D23").Copy($ws.Range("F1:F23"))

# ---------------------------------------------------------------------
# 2) Clear out the old D/E contents (they have been copied already) so
#    we can write the new NU_F / NU_V data into them.
# ---------------------------------------------------------------------
$ws.Range("D1:E23").ClearContents()

# ---------------------------------------------------------------------
# 3) New header row cells.
# ---------------------------------------------------------------------
$ws.Range("D1").Value = "NU_F"
$ws.Range("E1").Value = "NU_V"

# ---------------------------------------------------------------------
# 4) Fill NU_F / NU_V = 2.5 for the existing data rows (2-23), matching
#    the right-aligned numeric style already used in columns A-C/F.
# ---------------------------------------------------------------------
$dataRows = $ws.Range("D2:E23")
$dataRows.HorizontalAlignment = -4152
$dataRows.Value = 2.5

# ---------------------------------------------------------------------
# 5) Three brand-new result rows (24-26) from the latest grid-search run.
# ---------------------------------------------------------------------
$ws.Range("A24:F26").HorizontalAlignment = -4152

$ws.Range("A24").Value = 50
$ws.Range("B24").Value = 10
$ws.Range("C24").Value = 1
$ws.Range("D24").Value = 1.5
$ws.Range("E24").Value = 1.5
$ws.Range("F24").HorizontalAlignment = -4152
$ws.Range("F24").Formula = "'0.7757707945031895"
$ws.Range("F24").HorizontalAlignment = -4152

$ws.Range("A25").Value = 50
$ws.Range("B25").Value = 1
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 1.5
$ws.Range("E25").Value = 1.5

$ws.Range("A26").Value = 100
$ws.Range("B26").Value = 1
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 1.5
$ws.Range("E26").Value = 1.5

# F25 / F26 are plain, unstyled numeric results (no right-alignment like
# the rest of the table), so set them last, without ever touching their
# alignment.
$ws.Range("F25").Value = 0.77578327115172996
$ws.Range("F26").Value = 0.77578323012968098

# ---------------------------------------------------------------------
# 6) Column widths: shrink NU_F (D) and widen the relocated Score
#    column (F).
# ---------------------------------------------------------------------
$ws.Columns("D").ColumnWidth = 10.916666686534882
$ws.Columns("F").ColumnWidth = 17.583333333333332

# ---------------------------------------------------------------------
# 7) Selection moves to A27.
# ---------------------------------------------------------------------
$ws.Range("A27").Select()
